$d = $word.ActiveDocument
$wNS = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# 1. Title: "Defect Bug Report" -> "Bug Report"
$d.Content.Find.Execute("Defect Bug Report", $false, $false, $false, $false, $false, $true, 1, $false, "Bug Report", 2) | Out-Null

# 2. Submission date: "23/05/204" -> "24/05/204" (fix day typo)
$d.Content.Find.Execute("23/05/204", $false, $false, $false, $false, $false, $true, 1, $false, "24/05/204", 2) | Out-Null

# 3. Fix "Oganisation" -> "Organisation", add the missing period after "etc", and
#    place the _GoBack bookmark right after it (last-edit marker). This also
#    removes the now-stale spellcheck proofErr wrappers around "Oganisation"/"etc".
$rng = $d.Content
$rng.Find.Execute("Oganisation") | Out-Null
$xml = "<w:p xmlns:w='$wNS' w:rsidR='005C7DDA' w:rsidRDefault='005C7DDA'>" +
       "<w:pPr><w:spacing w:line='260' w:lineRule='atLeast'/><w:cnfStyle w:val='000000100000' w:firstRow='0' w:lastRow='0' w:firstColumn='0' w:lastColumn='0' w:oddVBand='0' w:evenVBand='0' w:oddHBand='1' w:evenHBand='0' w:firstRowFirstColumn='0' w:firstRowLastColumn='0' w:lastRowFirstColumn='0' w:lastRowLastColumn='0'/></w:pPr>" +
       "<w:r><w:t xml:space='preserve'>Fill all required information like Name, </w:t></w:r>" +
       "<w:r><w:t>Organisation</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'> Name, email etc.</w:t></w:r>" +
       "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
       "<w:r><w:t xml:space='preserve'> except Message</w:t></w:r>" +
       "</w:p>"
$rng.InsertXML($xml) | Out-Null

# 4. Fix "Alll" -> "All" (also removes the stale spellcheck proofErr wrapper)
$rng2 = $d.Content
$rng2.Find.Execute("Alll") | Out-Null
$xml2 = "<w:p xmlns:w='$wNS' w:rsidR='005C7DDA' w:rsidRDefault='005C7DDA'>" +
        "<w:pPr><w:spacing w:line='260' w:lineRule='atLeast'/><w:cnfStyle w:val='000000100000' w:firstRow='0' w:lastRow='0' w:firstColumn='0' w:lastColumn='0' w:oddVBand='0' w:evenVBand='0' w:oddHBand='1' w:evenHBand='0' w:firstRowFirstColumn='0' w:firstRowLastColumn='0' w:lastRowFirstColumn='0' w:lastRowLastColumn='0'/></w:pPr>" +
        "<w:r><w:t>All</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'> fields are filled in form</w:t></w:r>" +
        "</w:p>"
$rng2.InsertXML($xml2) | Out-Null
